# The source CSV/data feeding the "iris" sheet gained one extra leading
# data row (a placeholder "-" record) compared to the version already in
# the workbook, so the whole data block shifts down by one row and the
# new row is filled with "-" markers. The active sheet/selection also
# moved from "breast-cancer-wis" back to "iris".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("iris")
$ws2 = $wb.Worksheets.Item("breast-cancer-wis")

# Insert a new blank row at row 4, pushing the existing iris data (and the
# sheet's used-range dimension) down by one row.
$ws1.Rows("4:4").Insert() | Out-Null

# Stamp the newly inserted row with placeholder "-" values across C4:G4.
$ws1.Range("C4:G4").Value = "-"

# Record the new selection on the (now inactive) breast-cancer-wis sheet.
$ws2.Range("L21").Select() | Out-Null

# Switch the active sheet back to iris and set its new selection.
$ws1.Select() | Out-Null
$ws1.Range("G5").Select() | Out-Null
